# "Fruta / hortaliza, semanal" - weekly update: add a new week of price data
# for Hortaliza / Apio @ Terminal La Palmera de La Serena.
#
# A new pair of rows (Primera / Segunda calidad) is inserted at the top of
# the data block (row 401), pushing all existing rows down by two, and the
# sheet's used range grows from A1:R430 to A1:R432.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (row 401), shifting
# the existing rows (401-430) down to (403-432).
$ws.Rows.Item(401).Insert()
$ws.Rows.Item(401).Insert()

# --- Row 401: Apio, Americana (o), Primera ---
$ws.Range("A401").Value = 8
$ws.Range("B401").Value = "Terminal La Palmera de La Serena"
$ws.Range("C401").Value = "Coquimbo"
$ws.Range("D401").Value = 44746
$ws.Range("E401").Value = 4
$ws.Range("F401").Value = 100112017
$ws.Range("G401").Value = "Apio"
$ws.Range("H401").Value = "Americana (o)"
$ws.Range("I401").Value = "Primera"
$ws.Range("J401").Value = 2480
$ws.Range("K401").Value = 7800
$ws.Range("L401").Value = 8000
$ws.Range("M401").Value = 7900
$ws.Range("N401").Value = "$/docena de matas"
$ws.Range("O401").Value = "Provincia del Elquí"
$ws.Range("P401").Value = 1317
$ws.Range("Q401").Value = 6
$ws.Range("R401").Value = "Hortaliza"

# --- Row 402: Apio, Americana (o), Segunda ---
$ws.Range("A402").Value = 8
$ws.Range("B402").Value = "Terminal La Palmera de La Serena"
$ws.Range("C402").Value = "Coquimbo"
$ws.Range("D402").Value = 44746
$ws.Range("E402").Value = 4
$ws.Range("F402").Value = 100112017
$ws.Range("G402").Value = "Apio"
$ws.Range("H402").Value = "Americana (o)"
$ws.Range("I402").Value = "Segunda"
$ws.Range("J402").Value = 1500
$ws.Range("K402").Value = 6800
$ws.Range("L402").Value = 7000
$ws.Range("M402").Value = 6900
$ws.Range("N402").Value = "$/docena de matas"
$ws.Range("O402").Value = "Provincia del Elquí"
$ws.Range("P402").Value = 1150
$ws.Range("Q402").Value = 6
$ws.Range("R402").Value = "Hortaliza"
